$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9444041848182678
$ws.Range("B1").Value = 1.952302217483521
$ws.Range("C1").Value = 4.21762228012085
$ws.Range("D1").Value = 3.289098262786865
$ws.Range("E1").Value = 1.44475257396698
